$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "add dialog history panel" — the "@Shake" / "@Jump" animation-trigger marker
# lines move up to sit right after the line that precedes them, and the
# dialogue line that used to hold their slot moves down to follow them.
# Net effect inside the existing 4-row window (rows 8-11, column C):
#   C8: "呀吼! 怎麼了?"            -> "@Shake"
#   C9: "@Shake"                   -> "呀吼! 怎麼了?"
#   C10: "阿阿...阿魯大人讓我來找你" -> "@Jump"   (format swaps with C11 too)
#   C11: "@Jump"                   -> "阿阿...阿魯大人讓我來找你"

# Rows 8/9 only swap text; formatting (style) stays put on each row.
$c8 = $ws.Cells.Item(8, 3)
$c9 = $ws.Cells.Item(9, 3)
$tmp89 = $c8.Value
$c8.Value = $c9.Value
$c9.Value = $tmp89

# Rows 10/11 swap both text AND formatting (the marker's own style moves
# with it), so copy/paste formats along with the value swap.
$c10 = $ws.Cells.Item(10, 3)
$c11 = $ws.Cells.Item(11, 3)

$tmp1011 = $c10.Value

$c10.Copy()
$ws.Range("Z1").PasteSpecial(-4122) # xlPasteFormats : stash C10's format
$c11.Copy()
$c10.PasteSpecial(-4122)            # C10 <- C11's format
$ws.Range("Z1").Copy()
$c11.PasteSpecial(-4122)            # C11 <- stashed C10's format
$ws.Range("Z1").ClearContents()
$ws.Range("Z1").ClearFormats()

$c10.Value = $c11.Value
$c11.Value = $tmp1011

# Row heights follow the content: the short ASCII marker rows are 12.8pt,
# the Chinese dialogue rows auto-size to 14.15pt.
$ws.Rows(8).RowHeight = 12.8
$ws.Rows(9).RowHeight = 14.15
$ws.Rows(10).RowHeight = 12.8
$ws.Rows(11).RowHeight = 14.15

$ws.Range("C11").Select()
